# BOM-DM0260(Production_1V0).xlsx update
#
# 1) Remove the two unused blank worksheets ("Sheet2", "Sheet3") that were
#    left over from the Altium BOM export - only the BOM sheet remains.
# 2) Update the BOM line for connector J2 (row 9): the camera-module part
#    referenced by the connector's Footprint/DesignItemId/LibRef columns is
#    swapped from the IMX378 camera ("AC-PY004-IMX378" /
#    "AC_PY004-IMX378_CON") to the OV9282 camera ("AC-PY003-OV9282" /
#    "CMP-005-000036-1"), and the Name column is updated to match the new
#    connector part name ("AC_PY003-OV9282_CON").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM-DM0260(Production_1V0)")

# --- 1) Drop the empty placeholder sheets ---------------------------------
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Sheet2").Delete()
$null = $wb.Worksheets.Item("Sheet3").Delete()

# --- 2) Update row 9 (designator J2) --------------------------------------
$ws.Range("A9").Value = "AC_PY003-OV9282_CON"
$ws.Range("G9").Value = "AC-PY003-OV9282"
$ws.Range("H9").Value = "CMP-005-000036-1"
$ws.Range("I9").Value = "CMP-005-000036-1"

# Re-apply the row's original cell format (border + text style) to the
# cells we just rewrote, since setting .Value resets it to the workbook's
# default style. Pull the format from an untouched neighbour in the same
# row (column D) which already carries the correct style.
$ws.Range("D9").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("D9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
